# Auto-generated edit script: update cryptos list values (prices, volume%, and a couple of
# row-content swaps) to match the target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = '''69.536.69'
$ws.Range("E2").Value = '''  -0.15%  '

# Row 3
$ws.Range("D3").Value = '''3.674.12'
$ws.Range("E3").Value = '''  -0.97%  '

# Row 4
$ws.Range("E4").Value = '''  +0.02%  '

# Row 5
$ws.Range("D5").Value = '''618.28'
$ws.Range("E5").Value = '''  -8.14%  '

# Row 6
$ws.Range("D6").Value = '''159.53'
$ws.Range("E6").Value = '''  -1.47%  '

# Row 7
$ws.Range("E7").Value = '''  +0.04%  '

# Row 8
$ws.Range("E8").Value = '''  -0.60%  '

# Row 9
$ws.Range("E9").Value = '''  -1.75%  '

# Row 10
$ws.Range("D10").Value = '''7.19'
$ws.Range("E10").Value = '''  +1.16%  '

# Row 11
$ws.Range("D11").Value = '''0.440'
$ws.Range("E11").Value = '''  -0.87%  '

# Row 12
$ws.Range("E12").Value = '''  -2.73%  '

# Row 13
$ws.Range("D13").Value = '''4.293.34'
$ws.Range("E13").Value = '''  -1.08%  '

# Row 14
$ws.Range("D14").Value = '''32.48'
$ws.Range("E14").Value = '''  -1.24%  '

# Row 15
$ws.Range("D15").Value = '''3.684.99'
$ws.Range("E15").Value = '''  -1.16%  '

# Row 16
$ws.Range("D16").Value = '''69.592.00'
$ws.Range("E16").Value = '''  -0.13%  '

# Row 17
$ws.Range("D17").Value = '''0.118'
$ws.Range("E17").Value = '''  +0.66%  '

# Row 18
$ws.Range("E18").Value = '''  -0.17%  '

# Row 19
$ws.Range("E19").Value = '''  -2.60%  '

# Row 20
$ws.Range("D20").Value = '''10.29'
$ws.Range("E20").Value = '''  +4.86%  '

# Row 21
$ws.Range("D21").Value = '''469.72'
$ws.Range("E21").Value = '''  -0.91%  '

# Row 22
$ws.Range("D22").Value = '''0.649'
$ws.Range("E22").Value = '''  -0.84%  '

# Row 23
$ws.Range("D23").Value = '''79.45'
$ws.Range("E23").Value = '''  -1.30%  '

# Row 24
$ws.Range("D24").Value = '''3.820.25'
$ws.Range("E24").Value = '''  -1.06%  '

# Row 25
$ws.Range("E25").Value = '''  +0.06%  '

# Row 26
$ws.Range("D26").Value = '''0.0000123'
$ws.Range("E26").Value = '''  -4.04%  '

# Row 27
$ws.Range("E27").Value = '''  +0.43%  '

# Row 28
$ws.Range("E28").Value = '''  -4.67%  '

# Row 29
$ws.Range("D29").Value = '''2.60'
$ws.Range("E29").Value = '''  -3.47%  '

# Row 30
$ws.Range("E30").Value = '''  -4.16%  '

# Row 31
$ws.Range("E31").Value = '''  -0.12%  '

# Row 32
$ws.Range("D32").Value = '''1.97'
$ws.Range("E32").Value = '''  -2.27%  '

# Row 33
$ws.Range("D33").Value = '''26.61'
$ws.Range("E33").Value = '''  -1.22%  '

# Row 34
$ws.Range("E34").Value = '''  -3.08%  '

# Row 35
$ws.Range("E35").Value = '''  -3.19%  '

# Row 36
$ws.Range("D36").Value = '''3.673.63'
$ws.Range("E36").Value = '''  -0.73%  '

# Row 37
$ws.Range("E37").Value = '''  -3.28%  '

# Row 38
$ws.Range("E38").Value = '''  -0.03%  '

# Row 39
$ws.Range("D39").Value = '''178.43'
$ws.Range("E39").Value = '''  +2.51%  '

# Row 40
$ws.Range("E40").Value = '''  -0.15%  '

# Row 41
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '''2.22'
$ws.Range("E41").Value = '''  -1.58%  '

# Row 42
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '''5.76'
$ws.Range("E42").Value = '''  -5.78%  '

# Row 43
$ws.Range("D43").Value = '''0.0891'
$ws.Range("E43").Value = '''  -2.63%  '

# Row 44
$ws.Range("E44").Value = '''  -1.57%  '

# Row 45
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '''46.81'
$ws.Range("E45").Value = '''  -0.64%  '

# Row 46
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''29.25'
$ws.Range("E46").Value = '''  +5.54%  '

# Row 47
$ws.Range("D47").Value = '''2.71'
$ws.Range("E47").Value = '''  -2.39%  '

# Row 48
$ws.Range("E48").Value = '''  -0.36%  '

# Row 49
$ws.Range("D49").Value = '''0.000264'
$ws.Range("E49").Value = '''  -6.68%  '

# Row 50
$ws.Range("E50").Value = '''  -4.77%  '

# Row 51
$ws.Range("E51").Value = '''  -6.52%  '
